$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Site Agent/Technician"
$ws.Range("D2").Value = "KUROS CONSTRUCTION SOLUTIONS LTD."
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "29/08/2024"

$ws.Range("B3").Value = "Site Engineer"
$ws.Range("D3").Value = "KUROS CONSTRUCTION SOLUTIONS LTD."
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "29/08/2024"

$ws.Range("B4").Value = "Trainee Engineer"
$ws.Range("D4").Value = "KUROS CONSTRUCTION SOLUTIONS LTD."
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "29/08/2024"

$ws.Range("B5").Value = "Chief Supervisor"
$ws.Range("D5").Value = "KUROS CONSTRUCTION SOLUTIONS LTD."
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "29/08/2024"

$ws.Range("B6").Value = "Site Clerk"
$ws.Range("D6").Value = "KUROS CONSTRUCTION SOLUTIONS LTD."
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "29/08/2024"

$ws.Range("B7").Value = "M & E Coordinator"
$ws.Range("D7").Value = "KUROS CONSTRUCTION SOLUTIONS LTD."
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "29/08/2024"

$ws.Range("B8").Value = "Handyman, Building Maintenance"
$ws.Range("D8").Value = "KUROS CONSTRUCTION SOLUTIONS LTD."
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "29/08/2024"

$ws.Range("B9").Value = "Carpenter"
$ws.Range("D9").Value = "ECO DECK LTD"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "30/08/2024"

$ws.Range("B10").Value = "Administrative Secretary"
$ws.Range("D10").Value = "OPP CONTRACTING LTD"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "30/08/2024"

$ws.Range("B11").Value = "Administrative Officer"
$ws.Range("D11").Value = "OPP CONTRACTING LTD"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "30/08/2024"

$ws.Range("B12").Value = "General Worker"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "30/08/2024"

$ws.Range("B13").Value = "Plumber and Pipe Fitter, General"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "30/08/2024"

$ws.Range("B14").Value = "Electrician"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "30/08/2024"

$ws.Range("B15").Value = "Site Supervisor"
$ws.Range("D15").Value = "OPP CONTRACTING LTD"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "30/08/2024"

$ws.Range("B16").Value = "Draughtsperson"
$ws.Range("D16").Value = "OPP CONTRACTING LTD"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "30/08/2024"

$ws.Range("B17").Value = "Purchasing Clerk"

$ws.Range("B18").Value = "Gardener"
$ws.Range("D18").Value = "RENOVATIA LTD"

$ws.Range("B19").Value = "Driver, Truck/Goods Vehicle/Chauffeur Poid-Lourd"
$ws.Range("D19").Value = "S M S CONTRACTING LTD"

$ws.Range("D20").Value = "SOLVCON LIMITED"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "31/08/2024"

$ws.Range("B21").Value = "Accountant"
$ws.Range("D21").Value = "METASIGN COMPANY LTD"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "31/08/2024"

$ws.Range("B22").Value = "Welder"
$ws.Range("D22").Value = "ECO DECK LTD"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "31/08/2024"

$ws.Range("B23").Value = "Cabinet Maker/Menuisier"
$ws.Range("D23").Value = "RBRB CONSTRUCTION LTD"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "31/08/2024"

$ws.Range("B24").Value = "Painter"
$ws.Range("D24").Value = "RBRB CONSTRUCTION LTD"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = "31/08/2024"

$ws.Range("B25").Value = "Secretary"
$ws.Range("D25").Value = "PRO-DESIGN ENGINEERING CONSULTANTS LTD"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "31/08/2024"

$ws.Range("B26").Value = "Factory Operator"
$ws.Range("D26").Value = "ECO DECK LTD"
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "31/08/2024"

$ws.Range("B27").Value = "Storekeeper"
$ws.Range("D27").Value = "METASIGN COMPANY LTD"

$ws.Range("B28").Value = "Receptionist"
$ws.Range("D28").Value = "PRO-DESIGN ENGINEERING CONSULTANTS LTD"

$ws.Range("B29").Value = "Administrative Clerk"
$ws.Range("D29").Value = "METASIGN COMPANY LTD"

$ws.Range("B30").Value = "Accounts clerk"
$ws.Range("D30").Value = "TOPBUILDER CO LTD"

$ws.Range("B31").Value = "Procurement/Purchasing Officer"
$ws.Range("D31").Value = "TOPBUILDER CO LTD"

$ws.Range("B32").Value = "Electrical Technician"
$ws.Range("D32").Value = "PROELEC ELECTRICAL LTD"

$ws.Range("B33").Value = "Human Resource Officer"
$ws.Range("D33").Value = "GREEN SCAFF CO LTD"

$ws.Range("B34").Value = "Human Resource Assistant"
$ws.Range("D34").Value = "SAFETY CONSTRUCTION CO LTD"

$ws.Range("B35").Value = "Engineer, Building Construction"
$ws.Range("D35").Value = "JIANGXI CONSTRUCTION INTERNATIONAL ENGINEERING (MAURITIUS) CO LTD"

$ws.Range("B36").Value = "Technician, General"
$ws.Range("D36").Value = "JIANGXI CONSTRUCTION INTERNATIONAL ENGINEERING (MAURITIUS) CO LTD"

$ws.Range("B37").Value = "Civil Engineer"
$ws.Range("D37").Value = "METASIGN COMPANY LTD"

$ws.Range("B38").Value = "Site Supervisor"
$ws.Range("D38").Value = "METASIGN COMPANY LTD"

$ws.Range("B39").Value = "Electrician"
$ws.Range("D39").Value = "PROELEC ELECTRICAL LTD"
$ws.Range("F39").NumberFormat = "@"
$ws.Range("F39").Value = "02/09/2024"

$ws.Range("B40").Value = "Supervisor, Plumbing and Pipe Fitting"
$ws.Range("D40").Value = "AKGM CONTRACTING LTD"
$ws.Range("F40").NumberFormat = "@"
$ws.Range("F40").Value = "10/09/2024"

$ws.Range("B41").Value = "Plumber and Pipe Fitter, General"
$ws.Range("D41").Value = "AKGM CONTRACTING LTD"
$ws.Range("F41").NumberFormat = "@"
$ws.Range("F41").Value = "10/09/2024"

$ws.Range("B42").Value = "Human Resource Assistant"
$ws.Range("D42").Value = "AKGM CONTRACTING LTD"
$ws.Range("F42").NumberFormat = "@"
$ws.Range("F42").Value = "12/09/2024"

$ws.Range("B43").Value = "M&E Engineer"
$ws.Range("D43").Value = "CIVELMEC GROUP LTD"
$ws.Range("F43").NumberFormat = "@"
$ws.Range("F43").Value = "13/09/2024"

$ws.Range("B44").Value = "Cleaner, Building/Office"
$ws.Range("D44").Value = "ARWAN ENTERPRISE LTD"
$ws.Range("F44").NumberFormat = "@"
$ws.Range("F44").Value = "20/09/2024"

$ws.Range("B45").Value = "Accounts clerk"
$ws.Range("D45").Value = "KUROS CONSTRUCTION SOLUTIONS LTD."
$ws.Range("F45").NumberFormat = "@"
$ws.Range("F45").Value = "20/09/2024"

$ws.Range("B46").Value = "Store Clerk"
$ws.Range("D46").Value = "PROSEC LTD"
$ws.Range("F46").NumberFormat = "@"
$ws.Range("F46").Value = "30/09/2024"

$ws.Range("B47").Value = "Administrative Clerk"
$ws.Range("D47").Value = "PROSEC LTD"
